$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Dcn"
$ws.Range("C2").Value = "Tlr2"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 6.615074666666668
$ws.Range("H2").Value = 19.845224
$ws.Range("I2").Value = 0.0008916467884469992
$ws.Range("J2").Value = 0.0008916467884469989
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 57.24915866666667
$ws.Range("N2").Value = 171.747476
$ws.Range("O2").Value = 0.9704198736548433
$ws.Range("P2").Value = 0.9704198736548435
$ws.Range("Q2").Value = 378.7074591838472
$ws.Range("R2").Value = 3408.367132654625
$ws.Range("S2").Value = 0.0008652717637894837
$ws.Range("T2").Value = 0.0008652717637894837

$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Dcn"
$ws.Range("C3").Value = "Tlr2"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 6.615074666666668
$ws.Range("H3").Value = 19.845224
$ws.Range("I3").Value = 0.0008916467884469992
$ws.Range("J3").Value = 0.0008916467884469989
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.569166
$ws.Range("N3").Value = 4.707498
$ws.Range("O3").Value = 0.02659864191768634
$ws.Range("P3").Value = 0.02659864191768634
$ws.Range("Q3").Value = 10.38015025439467
$ws.Range("R3").Value = 93.42135228955202
$ws.Range("S3").Value = 0.00002371659364295675
$ws.Range("T3").Value = 0.00002371659364295675

$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Dcn"
$ws.Range("C4").Value = "Tlr2"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 6.615074666666668
$ws.Range("H4").Value = 19.845224
$ws.Range("I4").Value = 0.0008916467884469992
$ws.Range("J4").Value = 0.0008916467884469989
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.1758903333333333
$ws.Range("N4").Value = 0.527671
$ws.Range("O4").Value = 0.002981484427470275
$ws.Range("P4").Value = 0.002981484427470276
$ws.Range("Q4").Value = 1.163527688144889
$ws.Range("R4").Value = 10.471749193304
$ws.Range("S4").Value = 0.000002658431014558611
$ws.Range("T4").Value = 0.00000265843101455861

$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Dcn"
$ws.Range("C5").Value = "Tlr2"
$ws.Range("D5").Value = "FAPs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 7285.701009
$ws.Range("H5").Value = 21857.103027
$ws.Range("I5").Value = 0.9820406017477925
$ws.Range("J5").Value = 0.9820406017477923
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 57.24915866666667
$ws.Range("N5").Value = 171.747476
$ws.Range("O5").Value = 0.9704198736548433
$ws.Range("P5").Value = 0.9704198736548435
$ws.Range("Q5").Value = 417100.2530621344
$ws.Range("R5").Value = 3753902.27755921
$ws.Range("S5").Value = 0.9529917166720191
$ws.Range("T5").Value = 0.9529917166720191

$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Dcn"
$ws.Range("C6").Value = "Tlr2"
$ws.Range("D6").Value = "sCs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 7285.701009
$ws.Range("H6").Value = 21857.103027
$ws.Range("I6").Value = 0.9820406017477925
$ws.Range("J6").Value = 0.9820406017477923
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.569166
$ws.Range("N6").Value = 4.707498
$ws.Range("O6").Value = 0.02659864191768634
$ws.Range("P6").Value = 0.02659864191768634
$ws.Range("Q6").Value = 11432.47430948849
$ws.Range("R6").Value = 102892.2687853965
$ws.Range("S6").Value = 0.02612094631451875
$ws.Range("T6").Value = 0.02612094631451874

$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Dcn"
$ws.Range("C7").Value = "Tlr2"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 7285.701009
$ws.Range("H7").Value = 21857.103027
$ws.Range("I7").Value = 0.9820406017477925
$ws.Range("J7").Value = 0.9820406017477923
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.1758903333333333
$ws.Range("N7").Value = 0.527671
$ws.Range("O7").Value = 0.002981484427470275
$ws.Range("P7").Value = 0.002981484427470276
$ws.Range("Q7").Value = 1281.484379040013
$ws.Range("R7").Value = 11533.35941136012
$ws.Range("S7").Value = 0.002927938761254582
$ws.Range("T7").Value = 0.002927938761254582

$ws.Range("A8").Value = "ECs"
$ws.Range("B8").Value = "Dcn"
$ws.Range("C8").Value = "Tlr2"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 126.6246363333333
$ws.Range("H8").Value = 379.873909
$ws.Range("I8").Value = 0.01706775146376063
$ws.Range("J8").Value = 0.01706775146376063
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 57.24915866666667
$ws.Range("N8").Value = 171.747476
$ws.Range("O8").Value = 0.9704198736548433
$ws.Range("P8").Value = 0.9704198736548435
$ws.Range("Q8").Value = 7249.153896555966
$ws.Range("R8").Value = 65242.38506900369
$ws.Range("S8").Value = 0.01656288521903486
$ws.Range("T8").Value = 0.01656288521903486

$ws.Range("A9").Value = "ECs"
$ws.Range("B9").Value = "Dcn"
$ws.Range("C9").Value = "Tlr2"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 126.6246363333333
$ws.Range("H9").Value = 379.873909
$ws.Range("I9").Value = 0.01706775146376063
$ws.Range("J9").Value = 0.01706775146376063
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.569166
$ws.Range("N9").Value = 4.707498
$ws.Range("O9").Value = 0.02659864191768634
$ws.Range("P9").Value = 0.02659864191768634
$ws.Range("Q9").Value = 198.6950740966313
$ws.Range("R9").Value = 1788.255666869682
$ws.Range("S9").Value = 0.0004539790095246358
$ws.Range("T9").Value = 0.0004539790095246358

$ws.Range("A10").Value = "ECs"
$ws.Range("B10").Value = "Dcn"
$ws.Range("C10").Value = "Tlr2"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 126.6246363333333
$ws.Range("H10").Value = 379.873909
$ws.Range("I10").Value = 0.01706775146376063
$ws.Range("J10").Value = 0.01706775146376063
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.1758903333333333
$ws.Range("N10").Value = 0.527671
$ws.Range("O10").Value = 0.002981484427470275
$ws.Range("P10").Value = 0.002981484427470276
$ws.Range("Q10").Value = 22.27204949288211
$ws.Range("R10").Value = 200.448445435939
$ws.Range("S10").Value = 0.00005088723520113531
$ws.Range("T10").Value = 0.00005088723520113531
